$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new paragraph right after the "Задание: ... продолжительность."
#    paragraph, containing the description of the produced class diagram.
# ---------------------------------------------------------------------

$taskIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Подсчитать продолжительность") {
        $taskIndex = $i
        break
    }
    $i = $i + 1
}

$taskPara = $d.Paragraphs($taskIndex)
$taskPara.Range.InsertParagraphAfter()

$diagramPara = $d.Paragraphs($taskIndex + 1)
$diagramPara.Range.Text = "В ходе данной работы была разработана диаграмма классов, отражающая процесс звукозаписи. Основные элементы – классы, предназначенные для конкретного вида записей – песен, альбомов и сборников. Диаграмма представлена на рис.1."

# ---------------------------------------------------------------------
# 2) Insert a new paragraph at the very end of the document (after the
#    "Рисунок 1 – ..." caption), containing the "Вывод: ..." text. The
#    "UML" run must carry the en-US language mark, like elsewhere in the
#    document. Range.LanguageID in this runtime mis-scopes to the whole
#    paragraph, so the paragraph content is built from a raw OOXML
#    fragment via Range.InsertXML, which lets every run keep its own
#    run properties exactly as intended.
# ---------------------------------------------------------------------

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$conclusionPara = $d.Paragraphs($d.Paragraphs.Count)

# Seed paragraph/run formatting (justify + Times New Roman 14pt) so the
# paragraph mark itself carries the right rPr even though InsertXML below
# only rewrites the run content, not the paragraph mark.
$conclusionPara.Range.ParagraphFormat.Alignment = 3
$conclusionPara.Range.Font.Name = "Times New Roman"
$conclusionPara.Range.Font.Size = 14
$conclusionPara.Range.Font.SizeBi = 14
$conclusionPara.Range.Text = "x"

$contentRange = $d.Range($conclusionPara.Range.Start, $conclusionPara.Range.End - 1)

$conclusionXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Вывод:</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> По итогам данной работы была выполнена цель работы, заключающаяся в создании диаграммы классов </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>UML</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>по заданию.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$contentRange.InsertXML($conclusionXml)

Write-Output "done"
